$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lista de Itens de Trabalho")

# Update "Estado Atual" (D) percentages for existing rows
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 0.08
$ws.Range("D6").Value = 0.02

# Add two new work items
$ws.Range("A8").Value = "Criar tela de Login"
$ws.Range("B8").Value = "Alta"
$ws.Range("C8").Value = 40
$ws.Range("D8").Value = 0.8
$ws.Range("G8").Value = 72
$ws.Range("H8").Value = 0

$ws.Range("A9").Value = "Criar Localização de Veículos"
$ws.Range("B9").Value = "Alta"
$ws.Range("C9").Value = 60
$ws.Range("D9").Value = 0.2
$ws.Range("G9").Value = 120
$ws.Range("H9").Value = 0

# Match styles used by the other "Nome / Descrição" data rows (A5:A7 use style index 5)
$ws.Range("A8:A9").Style = $ws.Range("A5").Style
$ws.Range("D8:D9").Style = $ws.Range("D2").Style
$ws.Range("H8:H9").Style = $ws.Range("H2").Style

# Update selection to match the new active cell
$ws.Range("C10").Select()
